$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed by Excel as a number;
# force text storage via NumberFormat "@" then restore the Normal style so
# no visible/style diff is introduced (cell keeps default style index).

$ws.Range("D2").Value = "65.675.40"
$ws.Range("E2").Value = "  -0.95%  "
$ws.Range("D3").Value = "2.670.07"
$ws.Range("E3").Value = "  -0.83%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("E5").Value = "  -1.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.12"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.01%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.617"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.39%  "
$ws.Range("E9").Value = "  +4.06%  "
$ws.Range("E10").Value = "  -0.98%  "
$ws.Range("E11").Value = "  -2.69%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "29.46"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.73%  "
$ws.Range("E14").Value = "  -5.16%  "
$ws.Range("D15").Value = "3.151.66"
$ws.Range("D16").Value = "65.421.77"
$ws.Range("E16").Value = "  -1.05%  "
$ws.Range("D17").Value = "2.679.01"
$ws.Range("E17").Value = "  -0.56%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.88"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.14%  "
$ws.Range("E19").Value = "  -2.17%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.69"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.76%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "352.52"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.09%  "
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.83"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.45%  "
$ws.Range("E24").Value = "  +3.27%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.76"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.05%  "
$ws.Range("E26").Value = "  -1.98%  "
$ws.Range("E27").Value = "  -4.19%  "
$ws.Range("E28").Value = "  -5.36%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.11"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.46%  "
$ws.Range("E30").Value = "  +0.06%  "
$ws.Range("E31").Value = "  -3.06%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "531.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.79%  "
$ws.Range("E33").Value = "  -2.95%  "
$ws.Range("E34").Value = "  -1.75%  "
$ws.Range("E35").Value = "  +0.65%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.46"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.78%  "
$ws.Range("E38").Value = "  +0.01%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "158.82"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.82%  "
$ws.Range("E40").Value = "  -4.02%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "42.68"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.18%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "165.05"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.67%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.12"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.53%  "
$ws.Range("E45").Value = "  -1.22%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0611"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.61%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "23.14"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.70%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.644"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.14%  "
$ws.Range("E49").Value = "  -2.98%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.100"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.51%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "20.19"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.25%  "
